$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking strings
# (e.g. "64.925.61", "1.00", "169.59") are preserved exactly as text,
# matching the inline-string cell type used in the source data.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "64.925.61"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").Value = "3.542.35"
$ws.Range("E3").Value = "  +3.80%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "594.80"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").Value = "136.95"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("D7").Value = "3.538.43"
$ws.Range("E7").Value = "  +3.75%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  +2.66%  "
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("D11").Value = "6.97"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  +3.51%  "
$ws.Range("D13").Value = "4.144.73"
$ws.Range("E13").Value = "  +3.79%  "
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").Value = "27.20"
$ws.Range("E15").Value = "  +5.10%  "
$ws.Range("D16").Value = "3.544.69"
$ws.Range("E16").Value = "  +3.61%  "
$ws.Range("E17").Value = "  +1.27%  "
$ws.Range("D18").Value = "64.827.32"
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("D19").Value = "10.09"
$ws.Range("E19").Value = "  +7.50%  "
$ws.Range("D20").Value = "5.81"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("D21").Value = "14.24"
$ws.Range("E21").Value = "  +6.34%  "
$ws.Range("D22").Value = "388.71"
$ws.Range("E22").Value = "  +2.90%  "
$ws.Range("D23").Value = "0.575"
$ws.Range("E23").Value = "  +6.80%  "
$ws.Range("D24").Value = "3.686.84"
$ws.Range("E24").Value = "  +3.93%  "
$ws.Range("D25").Value = "73.95"
$ws.Range("E25").Value = "  +3.55%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +10.26%  "
$ws.Range("D28").Value = "7.69"
$ws.Range("E28").Value = "  +8.07%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +5.24%  "
$ws.Range("D31").Value = "8.22"
$ws.Range("E31").Value = "  +4.01%  "
$ws.Range("D32").Value = "3.554.13"
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("E33").Value = "  +19.03%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "23.84"
$ws.Range("E35").Value = "  +4.46%  "
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("D37").Value = "169.59"
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  +8.30%  "
$ws.Range("D39").Value = "6.86"
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("D40").Value = "4.97"
$ws.Range("E40").Value = "  +9.53%  "
$ws.Range("D41").Value = "0.0804"
$ws.Range("E41").Value = "  +7.09%  "
$ws.Range("E42").Value = "  +2.52%  "
$ws.Range("D43").Value = "26.67"
$ws.Range("E43").Value = "  +20.16%  "
$ws.Range("D44").Value = "42.82"
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  +10.50%  "
$ws.Range("D47").Value = "4.42"
$ws.Range("E47").Value = "  +4.31%  "
$ws.Range("E48").Value = "  +4.57%  "
$ws.Range("D49").Value = "2.438.87"
$ws.Range("E49").Value = "  +11.82%  "
$ws.Range("D50").Value = "6.87"
$ws.Range("E50").Value = "  +7.12%  "
$ws.Range("D51").Value = "303.29"
$ws.Range("E51").Value = "  +11.98%  "
